$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Probing the Enigma of Quantum Entanglement" "Delving into the Marvelous World of Chemistry: A Journey of Discovery and Transformation"

# --- Author name: "Rosalyn Barker" -> "Dr" + "." + " Susan Stevens" (3 runs) ---
$d.Content.Find.Execute("Rosalyn Barker", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Content
$rng.Start = $d.Content.Find.Parent.Start
$authorRng = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.Start)
# Locate via Find again to get a stable range object
$findRng = $d.Content
$findRng.Find.Execute("Rosalyn Barker", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nameStart = $findRng.Start
$nameEnd = $findRng.End
$nameFullRng = $d.Range($nameStart, $nameEnd)
$nameFullRng.Text = "Dr"
$afterDr = $nameStart + 2
$dotRng = $d.Range($afterDr, $afterDr)
$dotRng.InsertAfter(".")
$afterDot = $afterDr + 1
$restRng = $d.Range($afterDot, $afterDot)
$restRng.InsertAfter(" Susan Stevens")

# --- Email ---
Replace-Text "r_barker@amail" "susanstevensphds@gmail"
Replace-Text "science" "com"

# --- Body paragraph 1 (first set of sentences) ---
Replace-Text "Throughout the vast expanse of scientific inquiry, the phenomenon of quantum entanglement stands as an enigmatic paradox that challenges our understanding of the universe" "In the vast realm of science, chemistry stands as a captivating force, unveiling the mysteries of matter and its remarkable transformations"

Replace-Text " It is a realm where particles exhibit a profound interconnectedness, their destinies entwined in a manner that defies classical notions of locality and causality" " Embark with us on a captivating journey to unravel the enigmatic symphony of atoms and molecules, delving into the wonders of this dynamic field"

Replace-Text " This peculiar phenomenon has captivated the imagination of scientists for decades, sparking fervent debates and inspiring groundbreaking experiments that push the boundaries of human knowledge" " Chemistry permeates every aspect of our existence, from the air we breathe to the food we consume, orchestrating the symphony of life with unparalleled precision"

# Insert 4 new runs (., sentence, ., sentence) right after the run above and before the following "." run
$findRng2 = $d.Content
$findRng2.Find.Execute(" Chemistry permeates every aspect of our existence, from the air we breathe to the food we consume, orchestrating the symphony of life with unparalleled precision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPos = $findRng2.End
$insRng = $d.Range($insPos, $insPos)
$insRng.InsertAfter(".")
$insPos2 = $insPos + 1
$insRng2 = $d.Range($insPos2, $insPos2)
$insRng2.InsertAfter(" It unlocks the secrets hidden within the microscopic realm, revealing the essence of substances and their intricate interactions")
$insPos3 = $insPos2 + 129
$insRng3 = $d.Range($insPos3, $insPos3)
$insRng3.InsertAfter(".")
$insPos4 = $insPos3 + 1
$insRng4 = $d.Range($insPos4, $insPos4)
$insRng4.InsertAfter(" As we venture together into the world of chemistry, prepare to be captivated by the elegance and ingenuity with which matter transforms, revealing the profound interconnectedness of all things")

# --- Body paragraph 2 (after first double <br/>) ---
Replace-Text "In the microscopic realm of quantum physics, particles such as photons, electrons, and atoms can become entangled, forming a unified system that transcends the constraints of spatial separation" "Chemistry, the science of change, captivates with its breathtaking spectacles of transformation"

Replace-Text " The properties of these entangled particles become inextricably linked, exhibiting a remarkable correlation that extends across vast distances" " Witness the metamorphosis of substances, as they dance from one form to another, releasing vibrant hues and unleashing potent energies"

Replace-Text " This non-local connection between entangled particles has been experimentally confirmed through numerous experiments, including the groundbreaking work of John Clauser and Alain Aspect in the 1970s" " Whether it's the fiery combustion of fuels propelling rockets into the cosmos or the delicate synthesis of life-saving pharmaceuticals, chemistry orchestrates these transformations with exquisite precision"

Replace-Text " The results of these experiments have profoundly challenged our classical understanding of reality, leading to ongoing discussions about the fundamental nature of the universe" " By delving into its intricacies, we unravel the mysteries of matter, uncovering its hidden patterns and unlocking its boundless potential"

# --- Body paragraph 3 (after second double <br/>) ---
Replace-Text "The perplexing nature of quantum entanglement has profound implications for various fields of study, including quantum computing, cryptography, and potential applications in quantum teleportation" "In chemistry, we embark on a quest to understand not just what matter is but also how it behaves"

Replace-Text " The ability to harness and exploit the non-local correlations of entangled particles offers the potential to revolutionize communication, computation, and our understanding of the universe at its most fundamental level" " We decipher the enigmatic language of chemical reactions, unraveling the underlying principles that govern their interactions"

Replace-Text " As we delve deeper into the enigmatic realm of quantum entanglement, we continue to uncover its profound implications, pushing the boundaries of human knowledge and opening up new vistas of scientific exploration" " From the explosive reactions of fireworks illuminating the night sky to the controlled reactions within living cells fueling life's processes, chemistry unveils the symphony of molecular motion, energy transfer, and intricate equilibrium"

# --- Summary heading paragraph (body) ---
Replace-Text "Quantum entanglement, an enigmatic phenomenon that defies classical notions of locality and causality, has captivated the scientific community for decades" "Chemistry, a mesmerizing realm of discovery and transformation, captivates the imagination with its intricate symphony of matter"

Replace-Text " Entangled particles exhibit a profound interconnectedness, their properties becoming inextricably linked despite vast spatial " " It unveils the secrets of substances, orchestrating "

Replace-Text "separation" "breathtaking spectacles of change"

Replace-Text " Experiments like those conducted by John Clauser and Alain Aspect have confirmed the non-local nature of entanglement, challenging our understanding of reality" " Through chemistry, we decipher the enigmatic language of reactions, unraveling the underlying principles governing molecular interactions"

Replace-Text " This phenomenon has significant implications for various fields, including quantum computing, cryptography, and teleportation, offering the potential to revolutionize communication, computation, and our understanding of the universe" " From the fiery combustion of fuels to the delicate synthesis of pharmaceuticals, chemistry permeates every aspect of our existence, revealing the interconnectedness of all things"

Replace-Text " As we continue to explore the enigmatic realm of quantum entanglement, we uncover its profound implications, pushing the boundaries of human knowledge and opening up new avenues of scientific inquiry" " With each new discovery and each transformation, chemistry continues to redefine our understanding of the world around us, unlocking boundless possibilities for advancements in science, technology, and medicine"

# --- Add a new empty paragraph at the end of the document body (before sectPr) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "done"
